# Generate Report for Handoff
# Updates the localization-status report for the f68e5e0c file (row 3 in every
# sheet): the item is now ready for handoff again, with a fresh handoff
# timestamp and an error detail explaining the stale handback version.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e12470103b271a49d319fc31c15edc725f5f3f74/e2e/f68e5e0c-b2a4-47e5-8b0f-52d16917e7ba.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8e0bf7b9e8c96d2bd8dd5941d040df45ecac06b0/e2e/f68e5e0c-b2a4-47e5-8b0f-52d16917e7ba.md."

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-17 22:47:00"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-08-17 22:46:54"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.2

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-08-17 22:47:00"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.2
